$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 125 (this pushes the old
# rows 125 and 126 down to 127 and 128, and Excel auto-extends the
# dimension / keeps data intact).
$ws.Rows("125:126").Insert()

# --- New row 125: Ajo / Rosado / 1a nueva(o) ---
$ws.Range("A125").Value = 9
$ws.Range("B125").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C125").Value = "Metropolitana"
$ws.Range("D125").Value = 44509
$ws.Range("E125").Value = 13
$ws.Range("F125").Value = 100112003
$ws.Range("G125").Value = "Ajo"
$ws.Range("H125").Value = "Rosado"
$ws.Range("I125").Value = "1a nueva(o)"
$ws.Range("J125").Value = 52
$ws.Range("K125").Value = 3300
$ws.Range("L125").Value = 3600
$ws.Range("M125").Value = 3450
$ws.Range("N125").Value = "`$/paquete 20 unidades (volumen en unidades)"
$ws.Range("O125").Value = "Provincia de Talagante"
$ws.Range("P125").Value = 172
$ws.Range("Q125").Value = 20
$ws.Range("R125").Value = "Hortaliza"

# --- New row 126: Ajo / Rosado / Extra nueva (o) ---
$ws.Range("A126").Value = 9
$ws.Range("B126").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C126").Value = "Metropolitana"
$ws.Range("D126").Value = 44509
$ws.Range("E126").Value = 13
$ws.Range("F126").Value = 100112003
$ws.Range("G126").Value = "Ajo"
$ws.Range("H126").Value = "Rosado"
$ws.Range("I126").Value = "Extra nueva (o)"
$ws.Range("J126").Value = 79
$ws.Range("K126").Value = 3800
$ws.Range("L126").Value = 4000
$ws.Range("M126").Value = 3899
$ws.Range("N126").Value = "`$/paquete 20 unidades (volumen en unidades)"
$ws.Range("O126").Value = "Provincia de Talagante"
$ws.Range("P126").Value = 195
$ws.Range("Q126").Value = 20
$ws.Range("R126").Value = "Hortaliza"
